$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the displayed email text in column A (hyperlink targets stay as-is)
$ws.Range("A1").Value = "test1@izzum.com"
$ws.Range("A2").Value = "test2@gsc.com"
$ws.Range("A3").Value = "test3@sdfgdf.com"
$ws.Range("A4").Value = "test4@gsfd.com"

# Update the selected cell/range on the sheet
$ws.Range("E11").Select()
